# Insert a new data row at row 401 (pushes existing rows 401-424 down to 402-425)
# and populate it with the new weekly record for Vega Modelo de Temuco - Brócoli.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("401").Insert()

$ws.Cells.Item(401, 1).Value = 10
$ws.Cells.Item(401, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(401, 3).Value = "La Araucanía"
$ws.Cells.Item(401, 4).Value = 44706
$ws.Cells.Item(401, 5).Value = 9
$ws.Cells.Item(401, 6).Value = 100112023
$ws.Cells.Item(401, 7).Value = "Brócoli"
$ws.Cells.Item(401, 8).Value = "Sin especificar"
$ws.Cells.Item(401, 9).Value = "Primera"
$ws.Cells.Item(401, 10).Value = 2500
$ws.Cells.Item(401, 11).Value = 1100
$ws.Cells.Item(401, 12).Value = 1100
$ws.Cells.Item(401, 13).Value = 1100
$ws.Cells.Item(401, 14).Value = "$/unidad"
$ws.Cells.Item(401, 15).Value = "Región Metropolitana"
$ws.Cells.Item(401, 16).Value = 1100
$ws.Cells.Item(401, 17).Value = 1
$ws.Cells.Item(401, 18).Value = "Hortaliza"
